$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.273.03'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '3.387.93'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.16'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '178.84'
$ws.Range('E6').Value = '  +0.90%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  +8.27%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '48.33'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('E12').Value = '  +3.45%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '682.85'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').Value = '3.932.50'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').Value = '69.394.97'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.120'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.386.90'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.27'
$ws.Range('E20').Value = '  +0.82%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.908'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '17.11'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.36'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '101.16'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.71'
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '33.45'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.73'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E31').Value = '  +16.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '11.04'
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '549.21'
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '57.85'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = '3.601.74'
$ws.Range('E37').Value = '  -2.91%  '
$ws.Range('E38').Value = '  +3.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.39'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('E40').Value = '  +10.68%  '
$ws.Range('E41').Value = '  +4.92%  '
$ws.Range('E42').Value = '  +4.06%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('E44').Value = '  +3.37%  '
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '129.97'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('E51').Value = '  +1.75%  '
